$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) stays text so values like "67.778.39" or
# "0.0000179" are not coerced into numbers/scientific notation by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '67.778.39'
$ws.Range('E2').Value = '  +1.10%  '
$ws.Range('D3').Value = '3.508.72'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '598.96'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').Value = '180.81'
$ws.Range('E6').Value = '  +4.49%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.509.31'
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('D9').Value = '0.595'
$ws.Range('E9').Value = '  -0.78%  '
$ws.Range('E10').Value = '  +7.49%  '
$ws.Range('D11').Value = '7.14'
$ws.Range('E11').Value = '  -1.88%  '
$ws.Range('D12').Value = '0.435'
$ws.Range('E12').Value = '  +0.39%  '
$ws.Range('D13').Value = '4.113.07'
$ws.Range('E13').Value = '  +0.09%  '
$ws.Range('D14').Value = '32.46'
$ws.Range('E14').Value = '  +12.09%  '
$ws.Range('E15').Value = '  +1.48%  '
$ws.Range('D16').Value = '67.780.80'
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('D17').Value = '0.0000179'
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').Value = '3.509.89'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').Value = '14.43'
$ws.Range('E20').Value = '  +1.93%  '
$ws.Range('D21').Value = '397.60'
$ws.Range('E21').Value = '  +0.96%  '
$ws.Range('D22').Value = '7.94'
$ws.Range('E22').Value = '  -0.75%  '
$ws.Range('D23').Value = '73.43'
$ws.Range('E23').Value = '  +0.37%  '
$ws.Range('D24').Value = '0.540'
$ws.Range('E24').Value = '  +0.72%  '
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('E26').Value = '  +0.52%  '
$ws.Range('D27').Value = '0.0000122'
$ws.Range('E27').Value = '  +1.04%  '
$ws.Range('E28').Value = '  +3.26%  '
$ws.Range('D29').Value = '0.177'
$ws.Range('E29').Value = '  -2.41%  '
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('D31').Value = '6.23'
$ws.Range('E31').Value = '  -0.67%  '
$ws.Range('D32').Value = '1.44'
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('E33').Value = '  +1.50%  '
$ws.Range('D34').Value = '23.78'
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('D35').Value = '7.46'
$ws.Range('E35').Value = '  +1.49%  '
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('D37').Value = '1.61'
$ws.Range('E37').Value = '  -3.20%  '
$ws.Range('D38').Value = '163.97'
$ws.Range('E38').Value = '  +0.23%  '
$ws.Range('D39').Value = '0.879'
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('D40').Value = '1.91'
$ws.Range('E40').Value = '  +0.65%  '
$ws.Range('D41').Value = '2.75'
$ws.Range('E41').Value = '  +7.49%  '
$ws.Range('D42').Value = '6.94'
$ws.Range('E42').Value = '  -0.34%  '
$ws.Range('D43').Value = '4.69'
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('D44').Value = '2.865.67'
$ws.Range('E44').Value = '  +1.92%  '
$ws.Range('D45').Value = '26.24'
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('D46').Value = '0.0731'
$ws.Range('E46').Value = '  -2.04%  '
$ws.Range('D47').Value = '26.61'
$ws.Range('E47').Value = '  -2.81%  '
$ws.Range('D48').Value = '42.22'
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = '0.0302'
$ws.Range('E49').Value = '  +0.11%  '
$ws.Range('B50').Value = 'Bittensor'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D50').Value = '342.91'
$ws.Range('E50').Value = '  +1.96%  '
$ws.Range('E51').Value = '  -1.16%  '
